# GDE-7943: Added Facility Field Validation for ALERT_003
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALERT")

# Insert a fresh row 4 (pushing nothing down, sheet only has 3 rows of data)
# so the new row inherits the per-cell formatting of the row above it, just
# like Excel does for a row inserted/typed directly under existing data.
$ws.Rows.Item(4).Insert(-4121, 0)

# New row of data (row 4) mirroring the existing ALERT_002 row (row 3),
# but with a Facility Name/FCN column added to the validation string.
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "ALERT_003"
$ws.Cells.Item(4, 3).Value = $ws.Cells.Item(3, 3).Value()
$ws.Cells.Item(4, 4).Value = "|"
$ws.Cells.Item(4, 5).Value = "Facility_3"
$ws.Cells.Item(4, 6).Value = "Deal Name|Deal Tracking Number|Facility Name|Facility FCN|Alert Heading|Alert Content|User Name|Date Added / Amended"

# Update selection to reflect the new last-edited cell, as captured in the diff.
$ws.Range("F4").Select()
